$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(5, 2, 3, 3, 1, 4, 5, 8, 4, 7)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Range("C11").Select()
